$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 7000
$ws.Range("J51").Value = 6000
$ws.Range("L51").Value = 6000
$ws.Range("N51").Value = -6968
$ws.Range("H112").Value = 1039.9744
$ws.Range("I112").Value = 990
$ws.Range("J112").Value = 1041.2894
$ws.Range("K112").Value = 2970
$ws.Range("L112").Value = 3123.8682
$ws.Range("M112").Value = -1862
$ws.Range("N112").Value = -5339.8682
$ws.Range("H137").Value = 126762.03
$ws.Range("I137").Value = 167924.58
$ws.Range("J137").Value = 3274.375
$ws.Range("K137").Value = 503773.74
$ws.Range("L137").Value = 9823.125
$ws.Range("M137").Value = -501223.74
$ws.Range("N137").Value = -14923.125
$ws.Range("H138").Value = 2294.4521
$ws.Range("I138").Value = 1208.0741
$ws.Range("J138").Value = 2932.1086
$ws.Range("K138").Value = 3624.2223
$ws.Range("L138").Value = 8796.325800000001
$ws.Range("M138").Value = 1515.7777
$ws.Range("N138").Value = -19076.3258

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7207.035
$ws.Range("I32").Value = 5837.712
$ws.Range("J32").Value = 11963.632
$ws.Range("K32").Value = 5837.712
$ws.Range("L32").Value = 11963.632
$ws.Range("M32").Value = -5550.712
$ws.Range("N32").Value = -12537.632
$ws.Range("H61").Value = 2766
$ws.Range("I61").Value = 3161.5
$ws.Range("K61").Value = 3161.5
$ws.Range("M61").Value = -2949.5
$ws.Range("H74").Value = 32259680
$ws.Range("I74").Value = 41667360
$ws.Range("J74").Value = 4771.4287
$ws.Range("K74").Value = 41667360
$ws.Range("L74").Value = 4771.4287
$ws.Range("M74").Value = -41666486
$ws.Range("N74").Value = -6519.4287
$ws.Range("H77").Value = 32259680
$ws.Range("I77").Value = 41667360
$ws.Range("J77").Value = 4771.4287
$ws.Range("K77").Value = 208336800
$ws.Range("L77").Value = 23857.1435
$ws.Range("M77").Value = -208332432
$ws.Range("N77").Value = -32593.1435
$ws.Range("H102").Value = 1675.5714
$ws.Range("I102").Value = 1545.8
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 1545.8
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = 76.20000000000005
$ws.Range("N102").Value = -5244
$ws.Range("H132").Value = 14524.342
$ws.Range("I132").Value = 2142.4075
$ws.Range("K132").Value = 6427.2225
$ws.Range("M132").Value = -3897.2225
$ws.Range("H136").Value = 2766
$ws.Range("I136").Value = 3161.5
$ws.Range("K136").Value = 9484.5
$ws.Range("M136").Value = -6934.5

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H133").Value = 50000
$ws.Range("J133").Value = 50000
$ws.Range("L133").Value = 50000
$ws.Range("N133").Value = -60120

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 44073.418
$ws.Range("I58").Value = 4437.3335
$ws.Range("J58").Value = 57285.445
$ws.Range("K58").Value = 4437.3335
$ws.Range("L58").Value = 57285.445
$ws.Range("M58").Value = -4234.3335
$ws.Range("N58").Value = -57691.445
$ws.Range("H107").Value = 1626.0435
$ws.Range("I107").Value = 997.0909
$ws.Range("K107").Value = 997.0909
$ws.Range("M107").Value = 922.9091
$ws.Range("H136").Value = 44073.418
$ws.Range("I136").Value = 4437.3335
$ws.Range("J136").Value = 57285.445
$ws.Range("K136").Value = 13312.0005
$ws.Range("L136").Value = 171856.335
$ws.Range("M136").Value = -10762.0005
$ws.Range("N136").Value = -176956.335

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 816.5
$ws.Range("J97").Value = 879.8
$ws.Range("L97").Value = 2639.4
$ws.Range("N97").Value = -3631.4
$ws.Range("H98").Value = 1274.125
$ws.Range("I98").Value = 807.1667
$ws.Range("J98").Value = 2675
$ws.Range("K98").Value = 2421.5001
$ws.Range("L98").Value = 8025
$ws.Range("M98").Value = -923.5001000000002
$ws.Range("N98").Value = -11021
$ws.Range("H122").Value = 991.4375
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 991.4375
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 8922.9375
$ws.Range("N122").Value = -13822.9375
$ws.Range("H131").Value = 702.2
$ws.Range("J131").Value = 750.3333
$ws.Range("L131").Value = 2250.9999
$ws.Range("N131").Value = -12330.9999
$ws.Range("H137").Value = 18523522
$ws.Range("I137").Value = 1573.75
$ws.Range("J137").Value = 33341082
$ws.Range("K137").Value = 4721.25
$ws.Range("L137").Value = 100023246
$ws.Range("M137").Value = 378.75
$ws.Range("N137").Value = -100033446
$ws.Range("M122").ClearContents()

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 40346.8
$ws.Range("I132").Value = 8520.4
$ws.Range("J132").Value = 103999.6
$ws.Range("K132").Value = 25561.2
$ws.Range("L132").Value = 311998.8
$ws.Range("M132").Value = -23031.2
$ws.Range("N132").Value = -317058.8

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5453.846
$ws.Range("I7").Value = 5266.6665
$ws.Range("K7").Value = 5266.6665
$ws.Range("M7").Value = -5154.6665
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("H100").Value = 2657.8
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 2657.8
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 2657.8
$ws.Range("N100").Value = -3739.8
$ws.Range("H126").Value = 5453.846
$ws.Range("I126").Value = 5266.6665
$ws.Range("K126").Value = 15799.9995
$ws.Range("M126").Value = -13329.9995
$ws.Range("H136").Value = 2186.9524
$ws.Range("I136").Value = 2006.6316
$ws.Range("J136").Value = 3900
$ws.Range("K136").Value = 6019.8948
$ws.Range("L136").Value = 11700
$ws.Range("M136").Value = -3469.8948
$ws.Range("N136").Value = -16800
$ws.Range("M17").ClearContents()
$ws.Range("M100").ClearContents()

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2195.913
$ws.Range("I126").Value = 1789.7894
$ws.Range("K126").Value = 5369.3682
$ws.Range("M126").Value = -2899.3682
$ws.Range("H132").Value = 2050.611
$ws.Range("I132").Value = 1070.8
$ws.Range("J132").Value = 3275.375
$ws.Range("K132").Value = 3212.4
$ws.Range("L132").Value = 9826.125
$ws.Range("M132").Value = -682.3999999999996
$ws.Range("N132").Value = -14886.125
$ws.Range("H136").Value = 23462934
$ws.Range("I136").Value = 28674920
$ws.Range("J136").Value = 9000
$ws.Range("K136").Value = 86024760
$ws.Range("L136").Value = 27000
$ws.Range("M136").Value = -86022210
$ws.Range("N136").Value = -32100
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

